$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new departure entry as row 25 (NUMBER 24), mirroring the layout
# of the existing rows (e.g. row 9 / row 16 for the same flight W95178).
$row = 25

$ws.Cells.Item($row, 1).Value = 24                     # A - NUMBER
$ws.Cells.Item($row, 2).Value = "Monday, Jan 16"        # B - DATE
$ws.Cells.Item($row, 3).Value = "9:40 AM"               # C - TIME
$ws.Cells.Item($row, 4).Value = "W95178"                # D - FLIGHT
$ws.Cells.Item($row, 5).Value = "London"                # E - TO
$ws.Cells.Item($row, 6).Value = "(LTN)"                 # F - SHORT
$ws.Cells.Item($row, 7).Value = "Wizz Air "              # G - AIRLINE
$ws.Cells.Item($row, 8).Value = "A320"                  # H - MODEL
$ws.Cells.Item($row, 9).Value = "(G-WUKF)"               # I - AIRCFAT ID
$ws.Cells.Item($row, 10).Value = "9:38 AM"              # J - STATUS
$ws.Cells.Item($row, 11).Font.Bold = $false             # K - (blank, but present)
$ws.Cells.Item($row, 12).Value = "0 hours, -2 minutes"  # L - DIFFERENCE
$ws.Cells.Item($row, 13).Font.Bold = $false             # M - (blank, but present)
